$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-11 (employee_id, employee_name, department, absence_reason, absence_duration, absence_date, salary)
$data = @(
    @(13014, "Ana Liz Rezende", "Recursos Humanos", "Consulta medica", 3, 45089, 7640.04),
    @(7332, "Bernardo Camargo", "Marketing", "Consulta medica", 5, 45089, 5559.75),
    @(10498, "Heloísa Vargas", "Juridico", "Problemas pessoais", 3, 45080, 3296.63),
    @(27563, "Dra. Maria Cecília Silveira", "TI", "Problemas pessoais", 6, 45105, 4602.2),
    @(34182, "Dr. Gustavo Câmara", "Vendas", "Problemas pessoais", 8, 45098, 6261.49),
    @(62767, "Pietro Porto", "Vendas", "Problemas pessoais", 8, 45078, 3762.31),
    @(64609, "Sra. Sophie Borges", "Recursos Humanos", "Viagem de negocios", 2, 45095, 9289.27),
    @(25709, "Dra. Sara Montenegro", "Marketing", "Problemas pessoais", 1, 45104, 3969.7),
    @(91923, "Sra. Alice Novais", "Financeiro", "Doenca", 5, 45098, 3338.51),
    @(7336, "Leandro Cunha", "Vendas", "Consulta medica", 3, 45096, 7092.74)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $rowIndex++
}

$wb.Save()
